$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell that keeps the default (unstyled) format throughout the edit,
# used to restore style on cells where we temporarily force a text number format
# so that numeric-looking strings (e.g. "1.00") are not reinterpreted as numbers.
$defaultStyle = $ws.Range("D4").Style

$ws.Range("D2").Value = "68.751.97"
$ws.Range("E2").Value = "  +3.32%  "
$ws.Range("D3").Value = "2.556.98"
$ws.Range("E3").Value = "  +2.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.51"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +2.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.09"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  +3.28%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.534"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = "  +1.82%  "
$ws.Range("D9").Value = "2.555.18"
$ws.Range("E9").Value = "  +2.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.143"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  +3.13%  "
$ws.Range("E11").Value = "  +2.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.17"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.350"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.18"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = "  +2.08%  "
$ws.Range("D15").Value = "3.013.54"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000181"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = "  +2.69%  "
$ws.Range("D17").Value = "68.448.01"
$ws.Range("E17").Value = "  +3.23%  "
$ws.Range("D18").Value = "2.538.28"
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.07"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = "  +4.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.63"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  +3.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "368.16"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  +5.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.25"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = "  +1.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.77"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  +2.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.97"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("B25").Value = "Aptos"
$ws.Range("C25").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.37"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "  +3.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.31"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "  +2.31%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "2.675.47"
$ws.Range("E28").Value = "  +1.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.994"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0000100"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = "  +2.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "548.61"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  +3.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.34"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  +3.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.36"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "  +3.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.89"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "  +2.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.130"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.48"
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = "  +1.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.91"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.99"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = "  +2.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.70"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  +1.80%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.83"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  +3.00%  "
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.359"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  +1.42%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.27"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  +3.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.58"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  +2.03%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.569"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "  +2.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "148.72"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0283"
$ws.Range("E48").Value = "  +3.89%  "
$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.77"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = "  +2.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.73"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0760"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  +1.24%  "
